$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.471.55'
$ws.Range("E2").Value = '  -5.08%  '
$ws.Range("D3").Value = '2.434.12'
$ws.Range("E3").Value = '  -6.90%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '540.09'
$ws.Range("E5").Value = '  -6.01%  '
$ws.Range("D6").Value = '144.03'
$ws.Range("E6").Value = '  -7.92%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.11%  '
$ws.Range("D8").Value = '0.605'
$ws.Range("E8").Value = '  -3.03%  '
$ws.Range("D9").Value = '2.437.97'
$ws.Range("E9").Value = '  -6.70%  '
$ws.Range("E10").Value = '  -11.42%  '
$ws.Range("E11").Value = '  -2.13%  '
$ws.Range("D12").Value = '5.30'
$ws.Range("E12").Value = '  -9.33%  '
$ws.Range("E13").Value = '  -8.52%  '
$ws.Range("D14").Value = '25.60'
$ws.Range("E14").Value = '  -9.40%  '
$ws.Range("D15").Value = '2.876.74'
$ws.Range("E15").Value = '  -6.77%  '
$ws.Range("D16").Value = '60.449.58'
$ws.Range("E16").Value = '  -4.98%  '
$ws.Range("D17").Value = '0.0000161'
$ws.Range("E17").Value = '  -10.52%  '
$ws.Range("D18").Value = '2.446.18'
$ws.Range("E18").Value = '  -6.85%  '
$ws.Range("D19").Value = '10.97'
$ws.Range("E19").Value = '  -8.81%  '
$ws.Range("D20").Value = '6.86'
$ws.Range("E20").Value = '  -10.00%  '
$ws.Range("D21").Value = '4.13'
$ws.Range("E21").Value = '  -9.04%  '
$ws.Range("D22").Value = '313.77'
$ws.Range("E22").Value = '  -8.67%  '
$ws.Range("E23").Value = '  -0.01%  '
$ws.Range("D24").Value = '62.69'
$ws.Range("E24").Value = '  -7.23%  '
$ws.Range("E25").Value = '  -5.09%  '
$ws.Range("D26").Value = '2.597.91'
$ws.Range("E26").Value = '  -4.80%  '
$ws.Range("E27").Value = '  +0.34%  '
$ws.Range("D28").Value = '0.0₃0949'
$ws.Range("E28").Value = '  -12.86%  '
$ws.Range("D29").Value = '1.47'
$ws.Range("E29").Value = '  -7.02%  '
$ws.Range("D30").Value = '8.22'
$ws.Range("E30").Value = '  -10.43%  '
$ws.Range("D31").Value = '524.90'
$ws.Range("E31").Value = '  -11.17%  '
$ws.Range("D32").Value = '7.51'
$ws.Range("E32").Value = '  -5.07%  '
$ws.Range("D33").Value = '0.145'
$ws.Range("E33").Value = '  -10.10%  '
$ws.Range("E34").Value = '  -9.24%  '
$ws.Range("E35").Value = '  -11.21%  '
$ws.Range("D36").Value = '5.73'
$ws.Range("E36").Value = '  -12.93%  '
$ws.Range("D37").Value = '0.998'
$ws.Range("E37").Value = '  -0.12%  '
$ws.Range("D38").Value = '4.78'
$ws.Range("E38").Value = '  -10.98%  '
$ws.Range("E39").Value = '  -7.77%  '
$ws.Range("E40").Value = '  -8.11%  '
$ws.Range("D41").Value = '143.61'
$ws.Range("E41").Value = '  -6.84%  '
$ws.Range("D42").Value = '0.998'
$ws.Range("E42").Value = '  -0.16%  '
$ws.Range("E43").Value = '  -10.60%  '
$ws.Range("D44").Value = '39.87'
$ws.Range("E44").Value = '  -3.77%  '
$ws.Range("D45").Value = '2.26'
$ws.Range("E45").Value = '  -10.87%  '
$ws.Range("D46").Value = '144.85'
$ws.Range("E46").Value = '  -8.03%  '
$ws.Range("E47").Value = '  -9.49%  '
$ws.Range("D48").Value = '20.61'
$ws.Range("E48").Value = '  -13.19%  '
$ws.Range("D49").Value = '0.0524'
$ws.Range("E49").Value = '  -11.13%  '
$ws.Range("D50").Value = '0.0933'
$ws.Range("E50").Value = '  -6.91%  '
$ws.Range("D51").Value = '0.576'
$ws.Range("E51").Value = '  -8.70%  '
